# courseleafPatchControl.xlsx — "Added search to cim product defs"
#
# On the CIM sheet, insert two new process rows (a "searchCgi" /
# "command" pair, mirroring the existing ribbit.cgi pair) right after the
# existing ribbit/index.cgi pair, pushing the remaining rows down by two.
# Also update the active sheet/selection state to match (CIM becomes the
# active tab with a selection of B24; CAT keeps its old selection but is
# no longer the active tab).

$wb  = $excel.ActiveWorkbook
$cat = $wb.Worksheets.Item("CAT")
$cim = $wb.Worksheets.Item("CIM")

# --- CIM sheet: insert the two new rows ------------------------------------

# Make room: push rows 14.. down by two (two single-row inserts so the
# shift behaves the same as a 2-row block insert).
$cim.Rows.Item(14).Insert()
$cim.Rows.Item(14).Insert()

# Reuse the formatting of the existing cgi/command pair (rows 12:13 — the
# ribbit "/web/ribbit/index.cgi" block) for the two new rows.
$cim.Range("A12:D13").Copy()
$cim.Range("A14:D15").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# New "searchCgi" process pair.
$cim.Range("A14").Value = "searchCgi"
$cim.Range("B14").Value = "ribbit.cgi"
$cim.Range("C14").Value = "/web/search/index.cgi"

$cim.Range("A15").Value = "command"
$cim.Range("B15").Value = "chmod 750 ./web/search/index.cgi"
$cim.Range("C15").Value = "onChangeOnly"

# --- Selection / active-tab bookkeeping ------------------------------------

# Touch CAT first so its selection is recorded, then leave CIM active/
# selected last (matches the recorded workbook state).
$cat.Activate()
$x = $cat.Range("A22:XFD23").Select()

$cim.Activate()
$x = $cim.Range("B24").Select()
